# PowerShell Excel COM-interop script applying the crypto price / volume
# updates described by the target diff (commit message: "Updated cryptos
# list on Sun Sep 17 05:53:44 UTC 2023 with GitHub Actions").
#
# The workbook is already open as $excel.ActiveWorkbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text (e.g. "26.804.37", "0.500") that
# would otherwise be auto-coerced to a Number by Excel's smart-typing when
# assigned via .Value (since some of these look like valid numeric
# literals). Force the column to Text formatting first so every value is
# stored as a text string, matching the source data - then restore the
# "Normal" style so no stray number-format style lingers on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.804.37'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.645.18'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").Value = '216.78'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '0.500'
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = '0.0629'
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").Value = '19.19'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = '1.869.12'
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").Value = '1.658.68'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '4.17'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("E16").Value = '  -2.90%  '
$ws.Range("D17").Value = '26.803.97'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  -2.01%  '
$ws.Range("D19").Value = '214.01'
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").Value = '4.37'
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").Value = '2.38'
$ws.Range("E23").Value = '  +11.75%  '
$ws.Range("D24").Value = '9.38'
$ws.Range("E24").Value = '  -1.84%  '
$ws.Range("D25").Value = '145.09'
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").Value = '7.09'
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").Value = '15.69'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -1.13%  '
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").Value = '3.32'
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("D33").Value = '3.01'
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").Value = '1.289.40'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").Value = '0.0174'
$ws.Range("E37").Value = '  -5.51%  '
$ws.Range("D38").Value = '0.539'
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("D39").Value = '0.827'
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D43").Value = '5.36'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = '1.795.72'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").Value = '61.21'
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = '91.57'
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").Value = '0.0978'
$ws.Range("E51").Value = '  +0.01%  '

$priceRange.Style = "Normal"

$wb.Save()
